# Applies the "Add files via upload" edit:
#   - Removes the stray embedded "Picture 2" image from the title slide,
#     the "Synchronized parallelism" slide, and the "Sources" slide.
#   - Adds a new bibliography entry (Ахо, Лам, Сети, Ульман - Compilers)
#     as items 2 and 3 of the numbered list on the "Sources" slide,
#     right after the Rabaey reference.

$p = $ppt.ActivePresentation

# msoPicture constant used by Shape.Type
$msoPicture = 13

function Remove-PictureShapes($slide) {
    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Type -eq $msoPicture) {
            $shp.Delete()
        }
    }
}

# --- Slide 1 ("Параллельное программирование"): drop Picture 2 -------------
$slide1 = $p.Slides.Item(1)
Remove-PictureShapes $slide1

# --- Slide 5 ("Синхронизованный параллелизм"): drop Picture 2 --------------
$slide5 = $p.Slides.Item(5)
Remove-PictureShapes $slide5

# --- Slide 12 ("Источники"): drop Picture 2, extend the bibliography -------
$slide12 = $p.Slides.Item(12)
Remove-PictureShapes $slide12

# "Объект 2" is the bulleted placeholder holding the numbered reference list.
$refsShape = $slide12.Shapes.Item(2)
$refsRange = $refsShape.TextFrame.TextRange

# Append a new numbered paragraph with the "Ахо" reference, right after the
# existing Rabaey entry (new list formatting/numbering is inherited
# automatically from the paragraph being split, same as pressing Enter).
[void]$refsRange.InsertAfter("`rАхо")
$refsRange2 = $refsShape.TextFrame.TextRange
[void]$refsRange2.InsertAfter(" А.В., Лам М.С., Сети Р., Ульман Д.Д. Компиляторы: принципы, технологии и инструментарий. 2 изд. // Москва, Вильямс, 2018г. ")

# Add a trailing empty numbered paragraph (matches the source deck, which
# ends the list on a blank bulleted line). Insert a throwaway character to
# force the new paragraph mark, then delete it again so only the paragraph
# end-run properties remain.
$refsRange3 = $refsShape.TextFrame.TextRange
[void]$refsRange3.InsertAfter("`rX")
$refsRange4 = $refsShape.TextFrame.TextRange
$placeholder = $refsRange4.Characters($refsRange4.Length, 1)
$placeholder.Delete()
